# Update the localization-status workbook to reflect a new handoff report run.
# - Priority for the "372527d5-...", "63faa6c4-...", "bf75584e-...", "f118b84d-..."
#   source files changes from "low" to "ht" (rows 4-7) in both the zh-cn and
#   de-de target sheets.
# - The "Latest Handoff Datetime" for those same rows is refreshed to a new
#   generation timestamp for each locale.
# - Because that datetime value is also used by the Overview sheet's
#   "Latest HO Xliff Generate Date" column (it was a shared string in the
#   original file), update the Overview sheet too so the displayed values stay
#   consistent.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# New "Latest Handoff Datetime" values per locale (rows 4-7).
$newZhCnDate = "2016-08-15 10:32:30"
$newDeDeDate = "2016-08-15 10:32:35"

# --- zh-cn sheet ---
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = $newZhCnDate

# --- de-de sheet ---
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = $newDeDeDate

# --- Overview sheet (mirrors de-de's handoff datetime in this data set) ---
$wsOverview.Range("G4:G7").Value = $newDeDeDate
